# CronogramaFinal.xlsx edit
# The "Coleta de dados" block (2 rows) and the "Analise de dados" block
# (6 rows) swap positions inside the A8:F15 range - Coleta now comes
# first (rows 8-9), Analise follows (rows 10-15). Also updates the
# active sheet view (scroll position / selected cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New block contents (rows 8-15) -------------------------------------
# row -> (A, B, C(start-serial), D(end-serial), E, F)
$rows = @(
    ,@(8,  "Coleta ", "Coleta  de dados dos ensaios", 42856, 43040, "#DDEAD1", "black")
    ,@(9,  "Coleta ", "Coleta  de dados dos ensaios", 43101, 43160, "#DDEAD1", "black")
    ,@(10, "Análise ", " Análise de dados", 43009, 43069, "#4B6043", "white")
    ,@(11, "Análise ", " Análise de dados", 43101, 43434, "#4B6043", "white")
    ,@(12, "Análise ", " Análise de dados", 43466, 43799, "#4B6043", "white")
    ,@(13, "Análise ", " Análise de dados", 43831, 44165, "#4B6043", "white")
    ,@(14, "Análise ", " Análise de dados", 44197, 44499, "#4B6043", "white")
    ,@(15, "Análise ", " Análise de dados", 44562, 44681, "#4B6043", "white")
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
}

# --- Sheet view: scroll back to top-left and move selection to I10 ------
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("I10").Select()
